# Weekly update: insert a new pair of rows (Primera/Segunda) at the top of
# the data block (rows 178-179), shifting the existing historical rows
# down by two. The two oldest rows that fall off the bottom of the
# original range are appended as new rows 276-277.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 178; this shifts rows 178:275 down to
# 180:277, carrying their values/styles with them (including the trailing
# rows that now land on the newly extended range 276:277).
$ws.Rows("178:179").Insert()

# Populate the two new rows with this week's data. All of the
# "categorical" columns (A,B,C,E,F,G,H,I,N,O,Q,R) are constant across the
# whole table (I alternates Primera/Segunda by row parity), so copy them
# from the row directly below (the former row 178/179, now at 180/181) and
# just set the new date/price figures (D,J,K,L,M,P).

# Row 178 (Primera)
$ws.Range("A178").Value = $ws.Range("A180").Value()
$ws.Range("B178").Value = $ws.Range("B180").Value()
$ws.Range("C178").Value = $ws.Range("C180").Value()
$ws.Range("D178").Value = 44529
$ws.Range("E178").Value = $ws.Range("E180").Value()
$ws.Range("F178").Value = $ws.Range("F180").Value()
$ws.Range("G178").Value = $ws.Range("G180").Value()
$ws.Range("H178").Value = $ws.Range("H180").Value()
$ws.Range("I178").Value = $ws.Range("I180").Value()
$ws.Range("J178").Value = 2200
$ws.Range("K178").Value = 550
$ws.Range("L178").Value = 600
$ws.Range("M178").Value = 575
$ws.Range("N178").Value = $ws.Range("N180").Value()
$ws.Range("O178").Value = $ws.Range("O180").Value()
$ws.Range("P178").Value = 288
$ws.Range("Q178").Value = $ws.Range("Q180").Value()
$ws.Range("R178").Value = $ws.Range("R180").Value()

# Row 179 (Segunda)
$ws.Range("A179").Value = $ws.Range("A181").Value()
$ws.Range("B179").Value = $ws.Range("B181").Value()
$ws.Range("C179").Value = $ws.Range("C181").Value()
$ws.Range("D179").Value = 44529
$ws.Range("E179").Value = $ws.Range("E181").Value()
$ws.Range("F179").Value = $ws.Range("F181").Value()
$ws.Range("G179").Value = $ws.Range("G181").Value()
$ws.Range("H179").Value = $ws.Range("H181").Value()
$ws.Range("I179").Value = $ws.Range("I181").Value()
$ws.Range("J179").Value = 1440
$ws.Range("K179").Value = 450
$ws.Range("L179").Value = 500
$ws.Range("M179").Value = 475
$ws.Range("N179").Value = $ws.Range("N181").Value()
$ws.Range("O179").Value = $ws.Range("O181").Value()
$ws.Range("P179").Value = 238
$ws.Range("Q179").Value = $ws.Range("Q181").Value()
$ws.Range("R179").Value = $ws.Range("R181").Value()
